# Weekly data refresh: a new daily observation for "Berenjena" (Mercado
# Mayorista Lo Valledor de Santiago) is prepended to the data block, pushing
# every existing record (rows 279-316) down by one row and extending the
# used range from A1:R316 to A1:R317.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 279; Excel shifts rows 279:316 down to 280:317
# and inherits the row-above formatting (keeps the date-style on column D).
$ws.Rows.Item(279).Insert()

# Populate the newly inserted row with the latest observation.
$ws.Range("A279").Value = 6
$ws.Range("B279").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C279").Value = 'Metropolitana'
$ws.Range("D279").Value = 45077
$ws.Range("E279").Value = 13
$ws.Range("F279").Value = 100112001
$ws.Range("G279").Value = 'Berenjena'
$ws.Range("H279").Value = 'Sin especificar'
$ws.Range("I279").Value = 'Primera'
$ws.Range("J279").Value = 810
$ws.Range("K279").Value = 5000
$ws.Range("L279").Value = 6000
$ws.Range("M279").Value = 5519
$ws.Range("N279").Value = '$/caja 50 unidades'
$ws.Range("O279").Value = 'Región de Arica y Parinacota'
$ws.Range("P279").Value = 110
$ws.Range("Q279").Value = 50
$ws.Range("R279").Value = 'Hortaliza'
